# Applies the attendance-report sync update described by the commit:
#  - reorders the two "Recorded By" email lists (G2, G9)
#  - PARASITOLOGY / C1 / Session 1 (row 14) moves from Pending -> Recorded,
#    with a recorder e-mail and updated student count
#  - the Class Statistics (L6/L8/L9/L10) and Group Statistics (O15/Q15/R15/S15)
#    panels are recomputed to reflect the newly recorded session

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: reorder the recorder list (same people, new order) ---
$ws.Range("G2").Value = "System, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# --- Row 9: reorder the recorder list (same people, new order) ---
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# --- Class Statistics block ---
$ws.Range("L6").Value = 5              # Recorded Sessions: 4 -> 5
$ws.Range("L8").Value = 24             # Pending Sessions: 25 -> 24

$ws.Range("L9").NumberFormat = "@"     # keep as literal text, not a numeric percent
$ws.Range("L9").Value = "17.2%"        # Coverage %: 13.8% -> 17.2%

$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "21.6%"       # Average Attendance %: 23.8% -> 21.6%

# --- Row 14: PARASITOLOGY / C1 / Session 1 becomes Recorded ---
# Recolor A14:I14 from the "Pending" fill to the "Recorded" fill (match row 9's look).
$ws.Range("A14:I14").Interior.Color = $ws.Range("A9").Interior.Color()

$ws.Range("G14").Value = "esraa.sami@med.asu.edu.eg"
$ws.Range("H14").Value = "32/251"
$ws.Range("I14").Value = "Recorded"

# --- Group Statistics block (row 15: Year 2 / C1) ---
$ws.Range("O15").Value = 5             # Recorded: 4 -> 5
$ws.Range("Q15").Value = 24            # Pending: 25 -> 24

$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "17.2%"       # Coverage %: 13.8% -> 17.2%

$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "21.6%"       # Avg Attendance %: 23.8% -> 21.6%
